# cv123041a.xlsx - "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet "dados" originally contained two stray section-header rows
# (row 5 "situação do domicílio" and row 8 "grandes regiões e unidades da
# federação") that had no data of their own - they were leftovers from the
# source table and pushed all the real data rows down by two positions.
# The fix removes those two empty rows entirely (so every data row below
# shifts up and fills the gap), and renames the "unnamed: 1_level_1"
# column-2 header in row 2 to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header B2 from "unnamed: 1_level_1" to "total"
$ws.Range("B2").Value = "total"

# Remove the two empty/stray label rows. Delete the lower one (row 8)
# first so the row 5 index is still valid when we delete it next.
$ws.Rows.Item(8).EntireRow.Delete()
$ws.Rows.Item(5).EntireRow.Delete()
